# Branchements myRIO.xlsx - add AO (analog out) pin rows and GND row,
# reflowing the existing DIO/Demux/Encoder/PWM rows beneath them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Final B:C content for every row, 1-based row numbers -----------------
# Row 1 (header) is unchanged; rows 2-24 are the reflowed/extended table.
# New text is entered in the same order the original author typed it (GND
# row first, then the two A/AO rows, then the B/AO rows, then C/AO0) so the
# resulting shared-string table order lines up with the source workbook.
$ws.Cells.Item(11, 3).Value2 = "GND"
$ws.Cells.Item(11, 2).Value2 = "B/ Pin 6"
$ws.Cells.Item(3, 2).Value2  = "A/AO1 (Pin 4)"
$ws.Cells.Item(3, 3).Value2  = "VDD (cst =5)"
$ws.Cells.Item(2, 2).Value2  = "A/AO1 (Pin 2)"
$ws.Cells.Item(2, 3).Value2  = "VEE (cst = -10)"
$ws.Cells.Item(9, 2).Value2  = "B/AO0 (Pin 2)"
$ws.Cells.Item(10, 2).Value2 = "B/AO1 (Pin 4)"
$ws.Cells.Item(18, 2).Value2 = "C/AO0 (AO0)"

$rows = @(
    @(1,  "PIN",                    "Elément"),
    @(2,  "A/AO1 (Pin 2)",          "VEE (cst = -10)"),
    @(3,  "A/AO1 (Pin 4)",          "VDD (cst =5)"),
    @(4,  "A/DIO0 (Pin 11)",        "DIR"),
    @(5,  "A/DIO1 (Pin 13)",        "RST"),
    @(6,  "A/ENC (Pin 18 et 22)",   "Encodeur"),
    @(7,  "A/PWM (Pin 27)",         "PWM"),
    @(8,  "",                       ""),
    @(9,  "B/AO0 (Pin 2)",          "Demux: A"),
    @(10, "B/AO1 (Pin 4)",          "Demux: B"),
    @(11, "B/ Pin 6",               "GND"),
    @(12, "B/DIO0 (Pin 11)",        "ENA"),
    @(13, "B/DIO1 (Pin 13)",        "INH"),
    @(14, "B/DIO12 (Pin 22)",       "MS1"),
    @(15, "B/DIO14 (Pin 32)",       "MS2"),
    @(16, "B/DIO15 (Pin 34)",       "MS3"),
    @(17, "",                       ""),
    @(18, "C/AO0 (AO0)",            "Demux: C"),
    @(19, "C/DIO0 (DIO0)",          "SLP F"),
    @(20, "C/ DIO1 (DIO1)",         "SLP R"),
    @(21, "C/DIO2 (DIO2)",          "SLP B"),
    @(22, "C/DIO3 (DIO3)",          "SLP L"),
    @(23, "C/DIO4 (DIO4)",          "SLP D"),
    @(24, "C/DIO5 (DIO5)",          "SLP U")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
}

# --- Borders ----------------------------------------------------------
# Rows 2-21 already carried the thin-border style; rows 22-24 are brand
# new so they need the same thin border applied (color set before style
# so it reuses the existing border/style definitions instead of forking
# new ones).
$ws.Range("B22:C24").Borders.Color = 0
$ws.Range("B22:C24").Borders.LineStyle = 1

# --- Column C width -----------------------------------------------------
$ws.Columns("C").ColumnWidth = 12.3

# --- Selection ------------------------------------------------------------
$ws.Range("D17").Select() | Out-Null
